$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-22 Sunday" "2025-06-23 Monday"

Replace-Text "561÷8=70, 1" "818÷2=409, 0"
Replace-Text "676÷4=169, 0" "288÷6=48, 0"
Replace-Text "751÷6=125, 1" "828÷5=165, 3"
Replace-Text "933÷9=103, 6" "541÷8=67, 5"
Replace-Text "152÷2=76, 0" "230÷4=57, 2"

Replace-Text "753÷2=376, 1" "557÷6=92, 5"
Replace-Text "710÷6=118, 2" "884÷7=126, 2"
Replace-Text "958÷7=136, 6" "375÷3=125, 0"
Replace-Text "781÷3=260, 1" "891÷7=127, 2"
Replace-Text "540÷2=270, 0" "137÷7=19, 4"

Replace-Text "495÷7=70, 5" "793÷9=88, 1"
Replace-Text "836÷9=92, 8" "714÷5=142, 4"
Replace-Text "557÷8=69, 5" "674÷9=74, 8"
Replace-Text "285÷4=71, 1" "161÷2=80, 1"
Replace-Text "266÷9=29, 5" "704÷2=352, 0"

Replace-Text "633÷5=126, 3" "682÷9=75, 7"
Replace-Text "182÷9=20, 2" "185÷4=46, 1"
Replace-Text "584÷9=64, 8" "797÷4=199, 1"
Replace-Text "461÷8=57, 5" "318÷8=39, 6"
Replace-Text "485÷6=80, 5" "706÷9=78, 4"

Replace-Text "518÷4=129, 2" "702÷5=140, 2"
Replace-Text "553÷7=79, 0" "351÷6=58, 3"
Replace-Text "627÷3=209, 0" "784÷4=196, 0"
Replace-Text "371÷2=185, 1" "469÷9=52, 1"
Replace-Text "472÷9=52, 4" "562÷2=281, 0"
